$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Cells.Item(8, 1).Value = "Volume 32   Number  42"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- CompStat data table updates (rows 15-31) ---
# Row 15
$ws.Cells.Item(15, 3).Value = "'0"
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = -100
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 5
$ws.Cells.Item(15, 8).Value = -80
$ws.Cells.Item(15, 9).Value = 23
$ws.Cells.Item(15, 10).Value = 32
$ws.Cells.Item(15, 11).Value = -28.125
$ws.Cells.Item(15, 12).Value = 9.523809523809
$ws.Cells.Item(15, 13).Value = -11.538461538461
$ws.Cells.Item(15, 14).Value = -14.814814814814

# Row 16
$ws.Cells.Item(16, 3).Value = 8
$ws.Cells.Item(16, 4).Value = 12
$ws.Cells.Item(16, 5).Value = -33.333333333333
$ws.Cells.Item(16, 6).Value = 24
$ws.Cells.Item(16, 7).Value = 32
$ws.Cells.Item(16, 8).Value = -25
$ws.Cells.Item(16, 9).Value = 258
$ws.Cells.Item(16, 10).Value = 352
$ws.Cells.Item(16, 11).Value = -26.704545454545
$ws.Cells.Item(16, 12).Value = -7.857142857142
$ws.Cells.Item(16, 13).Value = -7.857142857142
$ws.Cells.Item(16, 14).Value = -75.66037735849

# Row 17
$ws.Cells.Item(17, 3).Value = 14
$ws.Cells.Item(17, 4).Value = 7
$ws.Cells.Item(17, 5).Value = 100
$ws.Cells.Item(17, 6).Value = 35
$ws.Cells.Item(17, 7).Value = 34
$ws.Cells.Item(17, 8).Value = 2.941176470588
$ws.Cells.Item(17, 9).Value = 435
$ws.Cells.Item(17, 10).Value = 521
$ws.Cells.Item(17, 11).Value = -16.506717850287
$ws.Cells.Item(17, 12).Value = 10.969387755102
$ws.Cells.Item(17, 13).Value = 62.31343283582
$ws.Cells.Item(17, 14).Value = 22.191011235955

# Row 18
$ws.Cells.Item(18, 3).Value = "'0"
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = -100
$ws.Cells.Item(18, 6).Value = 8
$ws.Cells.Item(18, 7).Value = 13
$ws.Cells.Item(18, 8).Value = -38.461538461538
$ws.Cells.Item(18, 9).Value = 125
$ws.Cells.Item(18, 10).Value = 183
$ws.Cells.Item(18, 11).Value = -31.693989071038
$ws.Cells.Item(18, 12).Value = 4.166666666666
$ws.Cells.Item(18, 13).Value = -47.478991596638
$ws.Cells.Item(18, 14).Value = -92.197253433208

# Row 19
$ws.Cells.Item(19, 3).Value = 13
$ws.Cells.Item(19, 4).Value = 19
$ws.Cells.Item(19, 5).Value = -31.578947368421
$ws.Cells.Item(19, 6).Value = 57
$ws.Cells.Item(19, 7).Value = 76
$ws.Cells.Item(19, 8).Value = -25
$ws.Cells.Item(19, 9).Value = 586
$ws.Cells.Item(19, 10).Value = 781
$ws.Cells.Item(19, 11).Value = -24.967989756722
$ws.Cells.Item(19, 12).Value = -13.696612665684
$ws.Cells.Item(19, 13).Value = 48.730964467005
$ws.Cells.Item(19, 14).Value = -50.127659574468

# Row 20
$ws.Cells.Item(20, 3).Value = 6
$ws.Cells.Item(20, 4).Value = 11
$ws.Cells.Item(20, 5).Value = -45.454545454545
$ws.Cells.Item(20, 6).Value = 23
$ws.Cells.Item(20, 7).Value = 37
$ws.Cells.Item(20, 8).Value = -37.837837837837
$ws.Cells.Item(20, 9).Value = 199
$ws.Cells.Item(20, 10).Value = 251
$ws.Cells.Item(20, 11).Value = -20.717131474103
$ws.Cells.Item(20, 12).Value = -27.372262773722
$ws.Cells.Item(20, 13).Value = 8.743169398907
$ws.Cells.Item(20, 14).Value = -88.999447208402

# Row 21
$ws.Cells.Item(21, 3).Value = 41
$ws.Cells.Item(21, 4).Value = 51
$ws.Cells.Item(21, 5).Value = -19.607843137254
$ws.Cells.Item(21, 6).Value = 148
$ws.Cells.Item(21, 7).Value = 197
$ws.Cells.Item(21, 8).Value = -24.8730964467
$ws.Cells.Item(21, 9).Value = 1632
$ws.Cells.Item(21, 10).Value = 2124
$ws.Cells.Item(21, 11).Value = -23.163841807909
$ws.Cells.Item(21, 12).Value = -7.692307692307
$ws.Cells.Item(21, 13).Value = 17.410071942446
$ws.Cells.Item(21, 14).Value = -72.998014559894

# Row 22
$ws.Cells.Item(22, 3).Value = "'0"
$ws.Cells.Item(22, 4).Value = "'0"
$ws.Cells.Item(22, 5).Value = "'***.*"
$ws.Cells.Item(22, 6).Value = 2
$ws.Cells.Item(22, 7).Value = 4
$ws.Cells.Item(22, 8).Value = -50
$ws.Cells.Item(22, 9).Value = 38
$ws.Cells.Item(22, 10).Value = 67
$ws.Cells.Item(22, 11).Value = -43.283582089552
$ws.Cells.Item(22, 12).Value = -47.222222222222
$ws.Cells.Item(22, 13).Value = 40.74074074074
$ws.Cells.Item(22, 14).Value = "'***.*"

# Row 24
$ws.Cells.Item(24, 3).Value = 28
$ws.Cells.Item(24, 4).Value = 20
$ws.Cells.Item(24, 5).Value = 40
$ws.Cells.Item(24, 6).Value = 113
$ws.Cells.Item(24, 7).Value = 117
$ws.Cells.Item(24, 8).Value = -3.418803418803
$ws.Cells.Item(24, 9).Value = 1157
$ws.Cells.Item(24, 10).Value = 1660
$ws.Cells.Item(24, 11).Value = -30.301204819277
$ws.Cells.Item(24, 12).Value = -27.278441231929
$ws.Cells.Item(24, 13).Value = 28.270509977827
$ws.Cells.Item(24, 14).Value = "'***.*"

# Row 25
$ws.Cells.Item(25, 3).Value = 11
$ws.Cells.Item(25, 4).Value = 8
$ws.Cells.Item(25, 5).Value = 37.5
$ws.Cells.Item(25, 6).Value = 36
$ws.Cells.Item(25, 7).Value = 48
$ws.Cells.Item(25, 8).Value = -25
$ws.Cells.Item(25, 9).Value = 384
$ws.Cells.Item(25, 10).Value = 867
$ws.Cells.Item(25, 11).Value = -55.709342560553
$ws.Cells.Item(25, 12).Value = -51.637279596977
$ws.Cells.Item(25, 13).Value = "'***.*"
$ws.Cells.Item(25, 14).Value = "'***.*"

# Row 26
$ws.Cells.Item(26, 3).Value = 13
$ws.Cells.Item(26, 4).Value = 24
$ws.Cells.Item(26, 5).Value = -45.833333333333
$ws.Cells.Item(26, 6).Value = 74
$ws.Cells.Item(26, 7).Value = 89
$ws.Cells.Item(26, 8).Value = -16.853932584269
$ws.Cells.Item(26, 9).Value = 846
$ws.Cells.Item(26, 10).Value = 1090
$ws.Cells.Item(26, 11).Value = -22.385321100917
$ws.Cells.Item(26, 12).Value = 9.302325581395
$ws.Cells.Item(26, 13).Value = 11.315789473684
$ws.Cells.Item(26, 14).Value = "'***.*"

# Row 27
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = -100
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = 7
$ws.Cells.Item(27, 8).Value = -85.714285714285
$ws.Cells.Item(27, 9).Value = 34
$ws.Cells.Item(27, 10).Value = 51
$ws.Cells.Item(27, 11).Value = -33.333333333333
$ws.Cells.Item(27, 12).Value = -15
$ws.Cells.Item(27, 13).Value = "'***.*"
$ws.Cells.Item(27, 14).Value = "'***.*"

# Row 28
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 5
$ws.Cells.Item(28, 7).Value = 9
$ws.Cells.Item(28, 8).Value = -44.444444444444
$ws.Cells.Item(28, 9).Value = 93
$ws.Cells.Item(28, 10).Value = 116
$ws.Cells.Item(28, 11).Value = -19.827586206896
$ws.Cells.Item(28, 12).Value = -23.770491803278
$ws.Cells.Item(28, 13).Value = "'***.*"
$ws.Cells.Item(28, 14).Value = "'***.*"

# Row 29
$ws.Cells.Item(29, 7).Value = "'0"
$ws.Cells.Item(29, 8).Value = "'***.*"
$ws.Cells.Item(29, 9).Value = 1
$ws.Cells.Item(29, 10).Value = 7
$ws.Cells.Item(29, 11).Value = -85.714285714285
$ws.Cells.Item(29, 12).Value = -80
$ws.Cells.Item(29, 13).Value = -83.333333333333
$ws.Cells.Item(29, 14).Value = -98.076923076923

# Row 30
$ws.Cells.Item(30, 7).Value = "'0"
$ws.Cells.Item(30, 8).Value = "'***.*"
$ws.Cells.Item(30, 9).Value = 1
$ws.Cells.Item(30, 10).Value = 5
$ws.Cells.Item(30, 11).Value = -80
$ws.Cells.Item(30, 12).Value = -75
$ws.Cells.Item(30, 13).Value = -80
$ws.Cells.Item(30, 14).Value = -97.916666666666

# Row 31
$ws.Cells.Item(31, 12).Value = -63.636363636363
